$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = 1
    4  = 5
    5  = -4
    6  = -5
    7  = -4
    8  = 6
    9  = 1
    10 = -3
    11 = 1
    12 = 5
    13 = 1
    14 = -1
    15 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
